# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the daily conversion note text ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$old = $ws1.Range("A1").Value()
$new = $old.Replace("1000 Bs = 1.77 = 6410.05 pesos", "1000 Bs = 1.83 = 6663.62 pesos")
$new = $new.Replace("6410.05 pesos = 1.75 = 916.47 Bs", "6663.62 pesos = 1.82 = 961.6 Bs")
$ws1.Range("A1").Value = $new

# --- tasas: update the rate table values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 547
$ws2.Range("O10").Value = 3645
$ws2.Range("N12").Value = 3661
$ws2.Range("O12").Value = 528.302
